$p = $ppt.ActivePresentation

# 1. Move slide 3 ("features & advantages of proposed method") to the end (position 8).
#    This is a reordering: everything that used to be in positions 4-8 shifts up to 3-7,
#    and the moved slide becomes the new position 8.
$p.Slides.Item(3).MoveTo(8)

# 2. Slide 2 ("survey of existing problems / challenges") - Content Placeholder body text edits.
$slide2 = $p.Slides.Item(2)
$content2 = $slide2.Shapes.Item(2)
$tr2 = $content2.TextFrame.TextRange

# Paragraph 2: "قالب های تعریف شده فقط بخشی از جمله را شامل میشوند"
#  -> "قالب های تعریف شده فقط بخشی از جمله را شامل میشوندمیشوند(قالب های تعریف شده محدود است ولی مال ما نامحدود)"
$tr2.Paragraphs(2, 1).Text = "قالب های تعریف شده فقط بخشی از جمله را شامل میشوندمیشوند(قالب های تعریف شده محدود است ولی مال ما نامحدود)"

# Paragraph 4: "استخراج روابط مفهومی از جملات ساده و در نظر نگرفتن جملات برزگ انجام میگرفته است"
#  -> "استخراج روابط مفهومی تنها از جملات ساده صورت میگرفته است"
$tr2.Paragraphs(4, 1).Text = "استخراج روابط مفهومی تنها از جملات ساده صورت میگرفته است"
